# Regenerate database with new stat order: PTS/REB/AST/STL/BLK
#
# The "stats" column (F) stores raw numbers in PTS/REB/AST/BLK/STL order,
# while the "bucket_desc" column (G) previously *labeled* them in a
# different order (PTS/AST/REB/BLK/STL). This edit re-labels G to match F's
# true order (PTS/REB/AST/.../...) and swaps the last two raw numbers in F
# (BLK and STL) so both columns present stats in PTS/REB/AST/STL/BLK order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 97
}

for ($row = 2; $row -le $lastRow; $row++) {

    # ---- Column F ("stats"): PTS/REB/AST/BLK/STL -> PTS/REB/AST/STL/BLK ----
    $fCell = $ws.Cells.Item($row, 6)
    $fVal = $fCell.Value2
    if ($fVal -ne $null -and $fVal -ne "") {
        $fParts = $fVal -split '/'
        if ($fParts.Length -eq 5) {
            $newF = $fParts[0] + "/" + $fParts[1] + "/" + $fParts[2] + "/" + $fParts[4] + "/" + $fParts[3]
            $fCell.Value2 = $newF
        }
    }

    # ---- Column G ("bucket_desc"): reorder labeled segments ----
    $gCell = $ws.Cells.Item($row, 7)
    $gVal = $gCell.Value2
    if ($gVal -ne $null -and $gVal -ne "") {
        $gParts = $gVal -split ' \| '
        if ($gParts.Length -eq 5) {
            $map = @{}
            foreach ($p in $gParts) {
                $kv = $p -split ' ', 2
                $map[$kv[0]] = $kv[1]
            }
            $newG = "PTS " + $map['PTS'] + " | REB " + $map['REB'] + " | AST " + $map['AST'] + " | STL " + $map['STL'] + " | BLK " + $map['BLK']
            $gCell.Value2 = $newG
        }
    }
}
